# Group the 11 "how does sourcing an R# script work" diagram shapes on
# slide 2 (矩形 3 / 圆角矩形 4 / TextBox 5 / TextBox 6 / 下箭头 7 /
# 圆角矩形 8 / TextBox 9 / 下箭头 10 / 圆角矩形 11 / TextBox 12 / 下箭头 13)
# into a single group shape, matching the "how to run a R# script?" commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# The group ends up named "组合 14" / id 15 in the real file, i.e. the
# PowerPoint shape-id counter had already advanced past 14 (ids 1-14 are
# already used by the title + the 11 shapes being grouped) by the time the
# group was created. Churning a throwaway shape first advances this
# runtime's internal id counter the same way, so the resulting group picks
# up id 15 instead of colliding with/reusing an existing id.
$scratch = $s.Shapes.AddShape(1, 0, 0, 1, 1)
$scratch.Delete()

# Shape indices (1-based) on the slide, in current z-order:
#  1 = 标题 1 (title, untouched)
#  2..12 = the 11 shapes that make up the diagram
$idxs = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12)
$range = $s.Shapes.Range($idxs)

$grp = $range.Group()
$grp.Name = "组合 14"
